$d = $word.ActiveDocument
$q  = [char]0x201C
$rq = [char]0x201D

# Helper: replace the full text of the first paragraph whose text matches
# $needle (substring match) with $newText. Using Range.Text directly (and
# not Find.Execute's ReplaceWith) avoids Word's smart-quote/apostrophe
# autocorrect mangling straight quotes in the replacement text.
# NB: iterate $d.Content.Paragraphs (recomputed from the Content range)
# rather than $d.Paragraphs, whose cached collection can go stale once a
# table has been structurally edited (e.g. Columns.Add/Rows.Add).
function Set-ParagraphText($needle, $newText) {
    foreach ($p in $d.Content.Paragraphs) {
        if ($p.Range.Text.Contains($needle)) {
            $p.Range.Text = $newText
            return $true
        }
    }
    return $false
}

# 1. Update the report date in the Heading1 paragraph.
Set-ParagraphText "April 04, 2024" "April 29, 2024" | Out-Null

# 2. Expand the italic query-instructions paragraph.
$newQuery = "Extract any quote that includes a national action or plan that addresses " + $q + "{variable_name}" + $rq + `
    " which we define as " + $q + "{variable_description}" + $rq + `
    ". Only include direct quotation with the corresponding page number(s) with a brief explanation of the context of this quote within the text. It is very important not to hallucinate."
Set-ParagraphText "Extract any quote that includes a national action or plan" $newQuery | Out-Null

# 3. First table ("query info"): add a 3rd column and a 2nd data row.
$t1 = $d.Tables.Item(1)
$t1.Columns.Add() | Out-Null
foreach ($col in $t1.Columns) {
    $col.Width = 144   # 2880 dxa == 144 pt
}

# Header row text updates.
$t1.Cell(1,2).Range.Text = "Variable description (optional)"
$hdr3 = $t1.Cell(1,3)
$hdr3.Range.InsertAfter("Context (optional)")
$hdr3r = $t1.Cell(1,3).Range
$hdr3Text = $d.Range($hdr3r.Start, $hdr3r.End - 1)
$hdr3Text.Font.Bold = 1

# First data row: "electricity grid" -> "SDG 1" plus a description.
$t1.Cell(2,1).Range.Text = "SDG 1"
$t1.Cell(2,2).Range.Text = "End poverty in all its forms everywhere."

# New second data row: "SDG 2".
$t1.Rows.Add() | Out-Null
$t1.Cell(3,1).Range.Text = "SDG 2"
$t1.Cell(3,2).Range.Text = "End hunger, achieve food security and improved nutrition and promote sustainable agriculture."

# 4. Source document heading.
Set-ParagraphText "ETH-008-lt-leds-2023.pdf" "Ghana_Ghana's Updated Nationally Determined Contribution to the UNFCCC_2021.pdf" | Out-Null

# 5. Second table (results): header rename + replace responses + add SDG 2 row.
$t2 = $d.Tables.Item(2)
$t2.Cell(1,2).Range.Text = "GPT Response"

$t2.Cell(2,1).Range.Text = "SDG 1"
$t2.Cell(2,2).Range.Text = "Build resilience and promote livelihood opportunities for the youth and women in climate- vulnerable Agriculture landscapes and food systems. [page(s) 26]"

$refRow = $t2.Rows.Item(3)
$t2.Rows.Add($refRow) | Out-Null
$t2.Cell(3,1).Range.Text = "SDG 2"
$t2.Cell(3,2).Range.Text = "Build resilience and promote livelihood opportunities for the youth and women in climate-vulnerable Agriculture landscapes and food systems. [page(s) 26]"

# 6. Footer processing-summary line.
Set-ParagraphText "1 documents (108 total pages) processed in 17.17 seconds" "1 documents (27 total pages) processed in 7.33 seconds" | Out-Null
